$d = $word.ActiveDocument

# The edit adds a single reviewer comment anchored on the whole first
# paragraph (from "Hola como estan, ..." through "... y aca quedamos."),
# exactly like the diff's new <w:commentRangeStart>/<w:commentRangeEnd>
# pair (w:id="0") around that paragraph's runs plus the new
# word/comments.xml part.

$p1 = $d.Paragraphs.Item(1)

$comment = $d.Comments.Add($p1.Range, "comment")
$comment.Author = "Obay Daba"
$comment.Initial = "od"
